# HyperParameter Tuned CartPole.xlsx - apply commit changes
# - Rename "Data" sheet to "CartPole"
# - Add a new "Trading" sheet with AnyTrading PPO results + bar chart
# - Update selections / active tab to match author's final view

$wb = $excel.ActiveWorkbook

# --- Rename existing sheet "Data" -> "CartPole" -----------------------
$cartpole = $wb.Worksheets.Item(1)
$cartpole.Name = "CartPole"

# --- Add the new "Trading" sheet after CartPole ------------------------
$trading = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $cartpole)
$trading.Name = "Trading"

# --- Populate the Trading sheet data -----------------------------------
$trading.Range("E3").Value = "Model Type"
$trading.Range("F3").Value = "AnyTrading Return (Training Data Only)"

$trading.Range("E4").Value = "Hyperparameter Tuned PPO (30000 timesteps)"
$trading.Range("F4").Value = 73.59

# Note: the "No Hyperparameter..." label is entered before the
# "Hyperparameter Tuned...(1000 timesteps)" label so that the shared
# string table ends up in the same order as the source workbook.
$trading.Range("E6").Value = "No Hyperparameter Tuning PPO (1000 timesteps)"
$trading.Range("F6").Value = -64.59

$trading.Range("E5").Value = "Hyperparameter Tuned PPO (1000 timesteps)"
$trading.Range("F5").Value = 73.59
$trading.Range("G5").Value = "*3x Improvement"

$trading.Range("E7").Value = "Random Action Sample"
$trading.Range("F7").Value = -96.59

$trading.Range("E3:F3").Font.Bold = $true

$trading.Columns.Item(5).ColumnWidth = 42.5546875
$trading.Columns.Item(6).ColumnWidth = 34.6640625

# --- Add the clustered-bar chart on the Trading sheet -------------------
$chartObj = $trading.ChartObjects().Add(0, 0, 400, 300)
$chart = $chartObj.Chart
$chart.ChartType = 51
$chart.SetSourceData($trading.Range("E4:F7"))
$chart.HasTitle = $false

$ser = $chart.SeriesCollection().Item(1)
$ser.Name = "=Trading!`$F`$3"

# --- Selections / active sheet to match the final saved view ------------
$cartpole.Range("C17").Select()
$trading.Activate()
$trading.Range("L19").Select()
